# Disabling TC23 & TC27 for Prod
#
# The sanity MasterExecutor sheet lists one automated test case per row
# (columns A:F). TC23_Verify_Footer and TC27_Price_Verification_on_CartPage
# are being disabled for the Prod run, so their rows are removed entirely
# (remaining rows shift up to close the gap), matching the commit intent
# "Disabling TC23 & TC27 for Prod".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the rows by their Testcase_number (column C) text instead of a
# hard-coded row index, so the edit still lands correctly even if the sheet
# layout shifts slightly.
$tc23Cell = $ws.Cells.Find("TC23_Verify_Footer")
$tc27Cell = $ws.Cells.Find("TC27_Price_Verification_on_CartPage")

$rowsToDelete = @()
if ($tc23Cell -ne $null) { $rowsToDelete += $tc23Cell.Row }
if ($tc27Cell -ne $null) { $rowsToDelete += $tc27Cell.Row }

# Delete from the bottom row up so earlier row numbers stay valid while we
# work through the list.
$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# Match the selection left behind by the author's editing session after the
# rows were removed.
$ws.Range("C2:C26").Select()
